$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New poll rows for the IFOP/OpinionWay rolling (2/15) and Harris/Cluster17 polls (2/15)

# Row 203: cluster17 / online / partially
$ws.Range("A203").Value = 100
$ws.Range("B203").Value = 2022
$ws.Range("C203").Value = 24
$ws.Range("D203").Value = 2
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = "cluster17"
$ws.Range("G203").Value = "online"
$ws.Range("H203").Value = "partially"
$ws.Range("I203").Value = 0
$ws.Range("J203").Value = 2085
$ws.Range("K203").Value = 0
$ws.Range("L203").Value = 1
$ws.Range("M203").Value = 0.5
$ws.Range("N203").Value = 0.5
$ws.Range("O203").Value = 13
$ws.Range("P203").Value = 4
$ws.Range("R203").Value = 5
$ws.Range("S203").Value = 2
$ws.Range("T203").Value = 23
$ws.Range("U203").Value = 15
$ws.Range("X203").Value = 0.5
$ws.Range("Y203").Value = 1
$ws.Range("Z203").Value = 16.5
$ws.Range("AA203").Value = 15
$ws.Range("AB203").Value = 0.5
$ws.Range("AD203").Value = 1
$ws.Range("AF203").Value = 2.5
$ws.Range("AH203").Value = 0

# Row 204: opinionway / online / partially
$ws.Range("A204").Value = 101
$ws.Range("B204").Value = 2022
$ws.Range("C204").Value = 25
$ws.Range("D204").Value = 2
$ws.Range("E204").Value = 14
$ws.Range("F204").Value = "opinionway"
$ws.Range("G204").Value = "online"
$ws.Range("H204").Value = "partially"
$ws.Range("I204").Value = 1
$ws.Range("J204").Value = 1000
$ws.Range("K204").Value = 1
$ws.Range("L204").Value = 1
$ws.Range("M204").Value = "T_1"
$ws.Range("N204").Value = "T_1"
$ws.Range("O204").Value = 10
$ws.Range("P204").Value = 5
$ws.Range("R204").Value = 5
$ws.Range("S204").Value = 3
$ws.Range("T204").Value = 25
$ws.Range("U204").Value = 16
$ws.Range("X204").Value = 2
$ws.Range("Y204").Value = 2
$ws.Range("Z204").Value = 15
$ws.Range("AA204").Value = 14
$ws.Range("AF204").Value = 3
$ws.Range("AH204").Value = 0

# Row 205: harris / online / included
$ws.Range("A205").Value = 102
$ws.Range("B205").Value = 2022
$ws.Range("C205").Value = 24
$ws.Range("D205").Value = 2
$ws.Range("E205").Value = 13
$ws.Range("F205").Value = "harris"
$ws.Range("G205").Value = "online"
$ws.Range("H205").Value = "included"
$ws.Range("I205").Value = 0
$ws.Range("J205").Value = 1700
$ws.Range("K205").Value = 1
$ws.Range("L205").Value = 1
$ws.Range("M205").Value = 0.5
$ws.Range("N205").Value = 0.5
$ws.Range("O205").Value = 10.5
$ws.Range("P205").Value = 3.5
$ws.Range("R205").Value = 5.5
$ws.Range("S205").Value = 2.5
$ws.Range("T205").Value = 25
$ws.Range("U205").Value = 14
$ws.Range("X205").Value = 1
$ws.Range("Y205").Value = 1
$ws.Range("Z205").Value = 17.5
$ws.Range("AA205").Value = 14.5
$ws.Range("AB205").Value = 0.5
$ws.Range("AF205").Value = 3.5
$ws.Range("AH205").Value = 0

# Row 206: ifop / online / included
$ws.Range("A206").Value = 103
$ws.Range("B206").Value = 2022
$ws.Range("C206").Value = 25
$ws.Range("D206").Value = 2
$ws.Range("E206").Value = 14
$ws.Range("F206").Value = "ifop"
$ws.Range("G206").Value = "online"
$ws.Range("H206").Value = "included"
$ws.Range("I206").Value = 1
$ws.Range("J206").Value = 1200
$ws.Range("K206").Value = 1
$ws.Range("L206").Value = 1
$ws.Range("M206").Value = 0.5
$ws.Range("N206").Value = 0.5
$ws.Range("O206").Value = 11.5
$ws.Range("P206").Value = 3
$ws.Range("R206").Value = 4.5
$ws.Range("S206").Value = 2
$ws.Range("T206").Value = 25.5
$ws.Range("U206").Value = 14.5
$ws.Range("X206").Value = 1
$ws.Range("Y206").Value = 1.5
$ws.Range("Z206").Value = 17.5
$ws.Range("AA206").Value = 15
$ws.Range("AF206").Value = 3
$ws.Range("AG206").Value = "T_0.5"
$ws.Range("AH206").Value = 0

# Update window/view scroll state to match where the author ended up after
# entering the new rows (scrolled back to column A, frozen pane near the
# new rows, active cell on the newly entered B206).
$win = $excel.ActiveWindow
$win.ScrollRow = 178
$win.ScrollColumn = 1
$win.Left = 22780
$win.Top = 8900
$ws.Range("B206").Select()
